# Se procesan de nuevo los datos con las nuevas dimensiones curadas
#
# The sheet stores, per data column (A..J), a small "metadata stack":
#   row1 = column name
#   row2 = semantic URI (iaest-measure:* / sdmx-dimension:* / null)
#   row3 = kind (medida / dim / null)
#   row4 = datatype / concept (xsd:int / skos:Concept / URI-* / null)
#   row5 = mapping workbook (only for dimension columns that used a lookup table)
#
# The curated re-processing turns the former "dimension" columns
# (inscripcion-provincia-nombre, sexo, inscripcion-comarca-nombre,
# edad-grandes-grupos) into "measure" columns: their semantic URI becomes
# iaest-measure:..., their kind becomes "medida" and their datatype becomes
# plain xsd:int (no more external URI/skos mapping workbook needed).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: semantic URI
$ws.Range("B2").Value = "iaest-measure:inscripcion-provincia-nombre"
$ws.Range("F2").Value = "iaest-measure:sexo"
$ws.Range("G2").Value = "iaest-measure:inscripcion-comarca-nombre"
$ws.Range("J2").Value = "iaest-measure:edad-grandes-grupos"

# Row 3: kind (dim -> medida)
$ws.Range("B3").Value = "medida"
$ws.Range("F3").Value = "medida"
$ws.Range("G3").Value = "medida"
$ws.Range("J3").Value = "medida"

# Row 4: datatype (URI-Provincia/URI-comarca/skos:Concept -> xsd:int)
$ws.Range("B4").Value = "xsd:int"
$ws.Range("F4").Value = "xsd:int"
$ws.Range("G4").Value = "xsd:int"
$ws.Range("J4").Value = "xsd:int"

# Row 5: the mapping-workbook references for sexo / edad-grandes-grupos no
# longer apply now that those columns are measures, not curated dimensions.
# Use Clear() (not ClearContents()) so the cell itself is dropped, matching
# row 5 having only the C5 (ano mapping) entry left.
$ws.Range("F5").Clear()
$ws.Range("J5").Clear()
